$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Alluvioni Emilia Romagna','FanPage','Facebook',86,14),
    @('Alluvioni Emilia Romagna','FanPage','Instagram',63,18),
    @('Alluvioni Emilia Romagna','FanPage','YouTube',67,31),
    @('Alluvioni Emilia Romagna','Il Corriere Della Sera','Facebook',72,35),
    @('Alluvioni Emilia Romagna','Il Corriere Della Sera','Instagram',65,6),
    @('Alluvioni Emilia Romagna','Il Corriere Della Sera','YouTube',82,19),
    @('Alluvioni Emilia Romagna','La Repubblica','Facebook',87,15),
    @('Alluvioni Emilia Romagna','La Repubblica','Instagram',72,13),
    @('Alluvioni Emilia Romagna','La Repubblica','YouTube',77,16),
    @('Attacco di Hamas a Israele','FanPage','Facebook',62,38),
    @('Attacco di Hamas a Israele','FanPage','Instagram',84,16),
    @('Attacco di Hamas a Israele','FanPage','YouTube',83,17),
    @('Attacco di Hamas a Israele','Il Corriere Della Sera','Facebook',82,18),
    @('Attacco di Hamas a Israele','Il Corriere Della Sera','Instagram',84,16),
    @('Attacco di Hamas a Israele','Il Corriere Della Sera','YouTube',81,19),
    @('Attacco di Hamas a Israele','La Repubblica','Facebook',86,13),
    @('Attacco di Hamas a Israele','La Repubblica','Instagram',85,8),
    @('Attacco di Hamas a Israele','La Repubblica','YouTube',78,17),
    @('Geolier e i fischi a Sanremo','FanPage','Facebook',33,67),
    @('Geolier e i fischi a Sanremo','FanPage','Instagram',75,25),
    @('Geolier e i fischi a Sanremo','FanPage','YouTube',80,20),
    @('Geolier e i fischi a Sanremo','Il Corriere Della Sera','Facebook',32,62),
    @('Geolier e i fischi a Sanremo','Il Corriere Della Sera','Instagram',78,22),
    @('Geolier e i fischi a Sanremo','Il Corriere Della Sera','YouTube',73,27),
    @('Geolier e i fischi a Sanremo','La Repubblica','Facebook',65,35),
    @('Geolier e i fischi a Sanremo','La Repubblica','Instagram',87,13),
    @('Geolier e i fischi a Sanremo','La Repubblica','YouTube',71,29),
    @('Il caso Chiara Ferragni e Balocco','FanPage','Facebook',84,10),
    @('Il caso Chiara Ferragni e Balocco','FanPage','Instagram',85,15),
    @('Il caso Chiara Ferragni e Balocco','FanPage','YouTube',86,14),
    @('Il caso Chiara Ferragni e Balocco','Il Corriere Della Sera','Facebook',49,48),
    @('Il caso Chiara Ferragni e Balocco','Il Corriere Della Sera','Instagram',92,8),
    @('Il caso Chiara Ferragni e Balocco','Il Corriere Della Sera','YouTube',71,19),
    @('Il caso Chiara Ferragni e Balocco','La Repubblica','Facebook',77,22),
    @('Il caso Chiara Ferragni e Balocco','La Repubblica','Instagram',89,11),
    @('Il caso Chiara Ferragni e Balocco','La Repubblica','YouTube',76,23),
    @('Incidente Youtubers','FanPage','Facebook',89,11),
    @('Incidente Youtubers','FanPage','Instagram',98,2),
    @('Incidente Youtubers','FanPage','YouTube',90,10),
    @('Incidente Youtubers','Il Corriere Della Sera','Facebook',94,6),
    @('Incidente Youtubers','Il Corriere Della Sera','Instagram',94,6),
    @('Incidente Youtubers','Il Corriere Della Sera','YouTube',94,6),
    @('Incidente Youtubers','La Repubblica','Facebook',86,14),
    @('Incidente Youtubers','La Repubblica','Instagram',91,9),
    @('Incidente Youtubers','La Repubblica','YouTube',93,7),
    @('Invasione Ucraina','FanPage','Facebook',94,6),
    @('Invasione Ucraina','FanPage','Instagram',91,9),
    @('Invasione Ucraina','FanPage','YouTube',79,21),
    @('Invasione Ucraina','Il Corriere Della Sera','Facebook',87,13),
    @('Invasione Ucraina','Il Corriere Della Sera','Instagram',89,11),
    @('Invasione Ucraina','Il Corriere Della Sera','YouTube',80,20),
    @('Invasione Ucraina','La Repubblica','Facebook',60,40),
    @('Invasione Ucraina','La Repubblica','Instagram',84,16),
    @('Invasione Ucraina','La Repubblica','YouTube',78,22),
    @('L''implosione del sottomarino Titan','FanPage','Facebook',88,12),
    @('L''implosione del sottomarino Titan','FanPage','Instagram',76,24),
    @('L''implosione del sottomarino Titan','FanPage','YouTube',77,23),
    @('L''implosione del sottomarino Titan','Il Corriere Della Sera','Facebook',89,11),
    @('L''implosione del sottomarino Titan','Il Corriere Della Sera','Instagram',83,17),
    @('L''implosione del sottomarino Titan','Il Corriere Della Sera','YouTube',84,16),
    @('L''implosione del sottomarino Titan','La Repubblica','Facebook',85,15),
    @('L''implosione del sottomarino Titan','La Repubblica','Instagram',87,13),
    @('L''implosione del sottomarino Titan','La Repubblica','YouTube',79,21),
    @('L''incoronazione di Re Carlo','FanPage','Facebook',45,55),
    @('L''incoronazione di Re Carlo','FanPage','Instagram',62,38),
    @('L''incoronazione di Re Carlo','FanPage','YouTube',67,33),
    @('L''incoronazione di Re Carlo','Il Corriere Della Sera','Facebook',69,31),
    @('L''incoronazione di Re Carlo','Il Corriere Della Sera','Instagram',78,22),
    @('L''incoronazione di Re Carlo','Il Corriere Della Sera','YouTube',75,25),
    @('L''incoronazione di Re Carlo','La Repubblica','Facebook',65,35),
    @('L''incoronazione di Re Carlo','La Repubblica','Instagram',72,28),
    @('L''incoronazione di Re Carlo','La Repubblica','YouTube',66,34),
    @('L''omicidio di Giulia Cecchettin','FanPage','Facebook',79,21),
    @('L''omicidio di Giulia Cecchettin','FanPage','Instagram',88,12),
    @('L''omicidio di Giulia Cecchettin','FanPage','YouTube',79,21),
    @('L''omicidio di Giulia Cecchettin','Il Corriere Della Sera','Facebook',91,9),
    @('L''omicidio di Giulia Cecchettin','Il Corriere Della Sera','Instagram',78,22),
    @('L''omicidio di Giulia Cecchettin','Il Corriere Della Sera','YouTube',88,12),
    @('L''omicidio di Giulia Cecchettin','La Repubblica','Facebook',87,13),
    @('L''omicidio di Giulia Cecchettin','La Repubblica','Instagram',87,13),
    @('L''omicidio di Giulia Cecchettin','La Repubblica','YouTube',82,18),
    @('La cattura di Matteo Messina Denaro','FanPage','Facebook',97,3),
    @('La cattura di Matteo Messina Denaro','FanPage','Instagram',81,19),
    @('La cattura di Matteo Messina Denaro','FanPage','YouTube',73,27),
    @('La cattura di Matteo Messina Denaro','Il Corriere Della Sera','Facebook',74,26),
    @('La cattura di Matteo Messina Denaro','Il Corriere Della Sera','Instagram',84,16),
    @('La cattura di Matteo Messina Denaro','Il Corriere Della Sera','YouTube',72,28),
    @('La cattura di Matteo Messina Denaro','La Repubblica','Facebook',84,15),
    @('La cattura di Matteo Messina Denaro','La Repubblica','Instagram',87,13),
    @('La cattura di Matteo Messina Denaro','La Repubblica','YouTube',78,22),
    @('La morte di Silvio Berlusconi','FanPage','Facebook',74,26),
    @('La morte di Silvio Berlusconi','FanPage','Instagram',75,25),
    @('La morte di Silvio Berlusconi','FanPage','YouTube',77,23),
    @('La morte di Silvio Berlusconi','Il Corriere Della Sera','Facebook',75,24),
    @('La morte di Silvio Berlusconi','Il Corriere Della Sera','Instagram',82,18),
    @('La morte di Silvio Berlusconi','Il Corriere Della Sera','YouTube',65,35),
    @('La morte di Silvio Berlusconi','Il Corrirere Della Sera','Facebook',1,0),
    @('La morte di Silvio Berlusconi','La Repubblica','Facebook',74,26),
    @('La morte di Silvio Berlusconi','La Repubblica','Instagram',81,19),
    @('La morte di Silvio Berlusconi','La Repubblica','YouTube',66,34),
    @('Strage di Cutro','FanPage','Facebook',93,7),
    @('Strage di Cutro','FanPage','Instagram',91,9),
    @('Strage di Cutro','FanPage','YouTube',87,13),
    @('Strage di Cutro','Il Corriere Della Sera','Facebook',90,10),
    @('Strage di Cutro','Il Corriere Della Sera','Instagram',92,8),
    @('Strage di Cutro','Il Corriere Della Sera','YouTube',88,12),
    @('Strage di Cutro','La Repubblica','Facebook',93,7),
    @('Strage di Cutro','La Repubblica','Instagram',79,21),
    @('Strage di Cutro','La Repubblica','YouTube',82,18)
)

$startRow = 3
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
}

Write-Output "Wrote $($data.Count) rows starting at row $startRow"
